$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new employee data column-by-column (A, B, C, D, E, F) across rows 6-21
# Column A: id
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20

# Column B: FullName
$ws.Range("B6").Value = 'Maya Lohar'
$ws.Range("B7").Value = 'Shyama rai'
$ws.Range("B8").Value = 'Gunjan Mehta'
$ws.Range("B9").Value = 'Vikram thakur'
$ws.Range("B10").Value = 'Abhimanyu singh'
$ws.Range("B11").Value = 'Suresh Kumar'
$ws.Range("B12").Value = 'Bharatram Patidar'
$ws.Range("B13").Value = 'Priya Sharma'
$ws.Range("B14").Value = 'Sapna Chouhan'
$ws.Range("B15").Value = 'Mohit Patidar'
$ws.Range("B16").Value = 'Kavita Laxmi'
$ws.Range("B17").Value = 'Swati Chouhan'
$ws.Range("B18").Value = 'Pradyuman Singh'
$ws.Range("B19").Value = 'Akshay Kumar'
$ws.Range("B20").Value = 'Twinkle Khanna'
$ws.Range("B21").Value = 'Sunita Patidar'

# Column C: DOB
$ws.Range("C6").Value = 34732
$ws.Range("C7").Value = 36883
$ws.Range("C8").Value = 34536
$ws.Range("C9").Value = 37061
$ws.Range("C10").Value = 37746
$ws.Range("C11").Value = 37735
$ws.Range("C12").Value = 35431
$ws.Range("C13").Value = 36681
$ws.Range("C14").Value = 36814
$ws.Range("C15").Value = 35912
$ws.Range("C16").Value = 37143
$ws.Range("C17").Value = 38112
$ws.Range("C18").Value = 37205
$ws.Range("C19").Value = 37094
$ws.Range("C20").Value = 36496
$ws.Range("C21").Value = 36506

# Column D: Salary
$ws.Range("D6").Value = 56000
$ws.Range("D7").Value = 23000
$ws.Range("D8").Value = 43000
$ws.Range("D9").Value = 56000
$ws.Range("D10").Value = 12000
$ws.Range("D11").Value = 49000
$ws.Range("D12").Value = 54000
$ws.Range("D13").Value = 23000
$ws.Range("D14").Value = 78000
$ws.Range("D15").Value = 67000
$ws.Range("D16").Value = 45000
$ws.Range("D17").Value = 55000
$ws.Range("D18").Value = 89000
$ws.Range("D19").Value = 36000
$ws.Range("D20").Value = 64000
$ws.Range("D21").Value = 56000

# Column E: Department
$ws.Range("E6").Value = 'Manager'
$ws.Range("E7").Value = 'IT'
$ws.Range("E8").Value = 'HR'
$ws.Range("E9").Value = 'IT'
$ws.Range("E10").Value = 'Finance'
$ws.Range("E11").Value = 'HR'
$ws.Range("E12").Value = 'HR'
$ws.Range("E13").Value = 'IT'
$ws.Range("E14").Value = 'Manager'
$ws.Range("E15").Value = 'IT'
$ws.Range("E16").Value = 'IT'
$ws.Range("E17").Value = 'Finance'
$ws.Range("E18").Value = 'IT'
$ws.Range("E19").Value = 'Health'
$ws.Range("E20").Value = 'IT'
$ws.Range("E21").Value = 'Manager'

# Column F: Age
$ws.Range("F6").Value = 29
$ws.Range("F7").Value = 24
$ws.Range("F8").Value = 30
$ws.Range("F9").Value = 23
$ws.Range("F10").Value = 21
$ws.Range("F11").Value = 21
$ws.Range("F12").Value = 27
$ws.Range("F13").Value = 24
$ws.Range("F14").Value = 24
$ws.Range("F15").Value = 26
$ws.Range("F16").Value = 23
$ws.Range("F17").Value = 19
$ws.Range("F18").Value = 23
$ws.Range("F19").Value = 22
$ws.Range("F20").Value = 25
$ws.Range("F21").Value = 25

# Apply existing date format (numFmtId 14) to C6:C20 by copying format from C2
$ws.Range("C2").Copy()
$ws.Range("C6:C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 21 date cell gets a new number format (numFmtId 16 => d-mmm)
$ws.Range("C21").NumberFormat = "d-mmm"

# Update selection to match the final cursor position
$ws.Range("D17").Select()
